# additions to feature tables
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 15 ("References" feature row): update status and add a comment
$ws.Range("B15").Value = "partially implemeted"
$ws.Range("C15").Value = "only ref. Checking, no “Go to definition”"

# Widen column C to fit the new comment text
# (engine quantizes ColumnWidth to ~1/6-character steps; 40.17 lands on the
#  stored width closest to the target 40.94)
$ws.Columns.Item(3).ColumnWidth = 40.17

# Reposition the active cell/selection as left by the author
$ws.Range("C5").Select() | Out-Null
